$wb = $excel.ActiveWorkbook

# 1. Update status text "Ready for handoff" -> "In Translation" on every sheet
#    (Overview: E2:F3, zh-cn: C2:C3, de-de: C2:C3 all reference this shared string)
foreach ($wsName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($wsName)
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# 2. Narrow the "Status" columns to match the new (shorter) text.
#    The runtime snaps ColumnWidth to an internal 1/6-character grid, so we
#    pick the input value that lands as close as possible to the target
#    stored width (~13.41 characters) -> 13.33 characters on the grid.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
